# Assignment 6 (Materials and Shading) -> Assignment 7 (Texture)
#
# The title paragraph's single run:
#     "Assignment 6 \u2013 Materials and Shading"
# becomes four separate runs reading "Assignment 7 \u2013 Texture":
#     "Assignment " | "7" | " \u2013 " | "Texture"

$d = $word.ActiveDocument
$dash = [char]0x2013   # en dash "\u2013"

# Locate the title text robustly via Find rather than hard-coded offsets.
$titleRange = $d.Content
$searchStr  = "Assignment 6 " + $dash + " Materials and Shading"
$found = $titleRange.Find.Execute($searchStr, $false, $false, $false, $false, $false, `
                                   $true, 1, $false, "", 0)

if ($found) {
    # Rebuild the paragraph's run content as four discrete <w:r> runs (rather than
    # one merged run) by injecting literal WordprocessingML for just this range;
    # the paragraph's own properties/identity (paraId, rsid, ...) are untouched
    # because the paragraph mark itself sits outside $titleRange.
    $newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
        '<w:r><w:t xml:space="preserve">Assignment </w:t></w:r>' +
        '<w:r><w:t>7</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> ' + $dash + ' </w:t></w:r>' +
        '<w:r><w:t>Texture</w:t></w:r>' +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $titleRange.InsertXML($newXml)
}
